$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern used throughout: the "Price" column (D) holds values that
# are stored as text in the source data (e.g. "235.45", "30.157.96" using
# dots as thousands separators). Setting .Value on a plain-numeric-looking
# string makes Excel auto-convert it to a number, so we force the cell to
# Text format first, assign the value, then drop the number format back to
# the workbook's default "Normal" style (so no stray formatting is left
# behind) while the stored value remains text.

# Row 19 / Row 20: the two coins (ShibaInu / Dai) swapped places in the
# ranking, so both rows' B/C/D/E contents need to move together.
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007518"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.36%  "

# Row 45 / Row 46: Quant and TheSandbox swapped places in the ranking.
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4137"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.56%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "

# Remaining per-row price / volume refreshes.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.157.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4721"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.25%  "

$ws.Range("E8").Value = "  +2.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06540"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07943"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.857.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.082"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6759"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.142.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.108.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.229"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.125"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.139"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.932"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("E29").Value = "  +1.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09874"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.465"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.288"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.996"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04672"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.120"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6975"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01871"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.606"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.316"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.920"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8364"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "943.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.153"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.976"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05650"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.38%  "
